$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$vB = New-Object 'object[,]' 24,1
$vB[0,0] = 1.019122443225342
$vB[1,0] = 0.9088306735465039
$vB[2,0] = 0.8409903731009649
$vB[3,0] = 0.8133163978711195
$vB[4,0] = 0.808719488668828
$vB[5,0] = 0.840617265133119
$vB[6,0] = 0.9811199667066148
$vB[7,0] = 1.255619955440977
$vB[8,0] = 1.456599289954056
$vB[9,0] = 1.54786533797369
$vB[10,0] = 1.582400779236764
$vB[11,0] = 1.57496409827246
$vB[12,0] = 1.550707102311947
$vB[13,0] = 1.535845674368204
$vB[14,0] = 1.450631442639235
$vB[15,0] = 1.398312777291096
$vB[16,0] = 1.368205477582308
$vB[17,0] = 1.358009149011593
$vB[18,0] = 1.403883753867945
$vB[19,0] = 1.557832664861564
$vB[20,0] = 1.658300383794312
$vB[21,0] = 1.604692945717375
$vB[22,0] = 1.401365204303488
$vB[23,0] = 1.181477577422982
$ws.Range("B2:B25").Value = $vB

$vC = New-Object 'object[,]' 24,1
$vC[0,0] = 0.2999826323890602
$vC[1,0] = 0.2828110259101493
$vC[2,0] = 0.2722194315562376
$vC[3,0] = 0.2678914742377572
$vC[4,0] = 0.2671721170021897
$vC[5,0] = 0.2721611105578461
$vC[6,0] = 0.2940720398292456
$vC[7,0] = 0.3366455048788453
$vC[8,0] = 0.3676715105706592
$vC[9,0] = 0.381728709278832
$vC[10,0] = 0.387043383963686
$vC[11,0] = 0.3858991561120604
$vC[12,0] = 0.3821661226221522
$vC[13,0] = 0.3798784175510548
$vC[14,0] = 0.3667516706273943
$vC[15,0] = 0.3586840820830162
$vC[16,0] = 0.3540385024360262
$vC[17,0] = 0.352464686039923
$vC[18,0] = 0.3595434432519369
$vC[19,0] = 0.383262837845848
$vC[20,0] = 0.3987152462821086
$vC[21,0] = 0.3904726505775216
$vC[22,0] = 0.3591549491428623
$vC[23,0] = 0.3251717320965213
$ws.Range("C2:C25").Value = $vC

$vD = New-Object 'object[,]' 24,1
$vD[0,0] = 0.07890202652542655
$vD[1,0] = 0.07152224508635641
$vD[2,0] = 0.0670248107635274
$vD[3,0] = 0.06520057099160681
$vD[4,0] = 0.06489817101962103
$vD[5,0] = 0.06700017399677449
$vD[6,0] = 0.07635046830709769
$vD[7,0] = 0.09495509251310352
$vD[8,0] = 0.1087901012453898
$vD[9,0] = 0.1151206752861071
$vD[10,0] = 0.1175232269217048
$vD[11,0] = 0.1170055589234806
$vD[12,0] = 0.1153182285050747
$vD[13,0] = 0.1142853797741736
$vD[14,0] = 0.1083771258719537
$vD[15,0] = 0.1047620598268679
$vD[16,0] = 0.1026862491108318
$vD[17,0] = 0.1019840133202621
$vD[18,0] = 0.105146529747131
$vD[19,0] = 0.1158136946077946
$vD[20,0] = 0.1228161926543407
$vD[21,0] = 0.1190760061200962
$vD[22,0] = 0.1049727029884622
$vD[23,0] = 0.08989303258667292
$ws.Range("D2:D25").Value = $vD

$vE = New-Object 'object[,]' 24,1
$vE[0,0] = 0.119416969029484
$vE[1,0] = 0.1208221098505957
$vE[2,0] = 0.1217347141096512
$vE[3,0] = 0.1221191548503998
$vE[4,0] = 0.122183749164873
$vE[5,0] = 0.1217398479906582
$vE[6,0] = 0.1198911253231126
$vE[7,0] = 0.1166605978202429
$vE[8,0] = 0.1145268486284626
$vE[9,0] = 0.1136079962466221
$vE[10,0] = 0.1132674853748142
$vE[11,0] = 0.1133404898138224
$vE[12,0] = 0.1135798331906821
$vE[13,0] = 0.1137274063386604
$vE[14,0] = 0.1145879392739745
$vE[15,0] = 0.1151291084232102
$vE[16,0] = 0.1154452509540742
$vE[17,0] = 0.1155531292165748
$vE[18,0] = 0.1150709954346273
$vE[19,0] = 0.113509330439853
$vE[20,0] = 0.1125320451296147
$vE[21,0] = 0.1130496772681364
$vE[22,0] = 0.115097252687949
$vE[23,0] = 0.1174923742754274
$ws.Range("E2:E25").Value = $vE

$vG = New-Object 'object[,]' 24,1
$vG[0,0] = 0.4792134293640444
$vG[1,0] = 0.4814597945129435
$vG[2,0] = 0.4832916227466981
$vG[3,0] = 0.4841516484494122
$vG[4,0] = 0.4843013051015532
$vG[5,0] = 0.483302761973448
$vG[6,0] = 0.4798938939242845
$vG[7,0] = 0.4768120221992831
$vG[8,0] = 0.4767617472620174
$vG[9,0] = 0.4772234163734339
$vG[10,0] = 0.4774682041188214
$vG[11,0] = 0.4774123684893112
$vG[12,0] = 0.4772421513930283
$vG[13,0] = 0.4771470083383917
$vG[14,0] = 0.4767413508957787
$vG[15,0] = 0.4766168005909179
$vG[16,0] = 0.4765907526397086
$vG[17,0] = 0.4765897545999422
$vG[18,0] = 0.4766253386907948
$vG[19,0] = 0.4772902472055023
$vG[20,0] = 0.4781327371164537
$vG[21,0] = 0.4776456639822584
$vG[22,0] = 0.4766213367286127
$vG[23,0] = 0.4772582131164711
$ws.Range("G2:G25").Value = $vG

$vH = New-Object 'object[,]' 24,1
$vH[0,0] = 0.6262227303661376
$vH[1,0] = 0.6318239643456351
$vH[2,0] = 0.6356278474994568
$vH[3,0] = 0.6372696449720578
$vH[4,0] = 0.6375478010100437
$vH[5,0] = 0.6356496181601656
$vH[6,0] = 0.628078327640381
$vH[7,0] = 0.6161260269462048
$vH[8,0] = 0.6091114739449353
$vH[9,0] = 0.6063044594087614
$vH[10,0] = 0.6052967620002363
$vH[11,0] = 0.6055113292904224
$vH[12,0] = 0.6062204478095481
$vH[13,0] = 0.6066620010256827
$vH[14,0] = 0.6093026489141948
$vH[15,0] = 0.6110209714555168
$vH[16,0] = 0.612045440612448
$vH[17,0] = 0.6123985131086442
$vH[18,0] = 0.610834312851992
$vH[19,0] = 0.6060106624429835
$vH[20,0] = 0.603180236422574
$vH[21,0] = 0.6046613993821524
$vH[22,0] = 0.6109185872562506
$vH[23,0] = 0.6190493118076716
$ws.Range("H2:H25").Value = $vH

$vI = New-Object 'object[,]' 24,1
$vI[0,0] = 0.5951951052498465
$vI[1,0] = 0.6051169995945287
$vI[2,0] = 0.611646899869811
$vI[3,0] = 0.614417838195811
$vI[4,0] = 0.6148845869980235
$vI[5,0] = 0.6116838247149499
$vI[6,0] = 0.5985252213346755
$vI[7,0] = 0.5762012167079398
$vI[8,0] = 0.5619292673487593
$vI[9,0] = 0.5559006902557222
$vI[10,0] = 0.553684672237889
$vI[11,0] = 0.5541589543192025
$vI[12,0] = 0.5557170358573202
$vI[13,0] = 0.5566801199611753
$vI[14,0] = 0.562332597264934
$vI[15,0] = 0.5659191278104601
$vI[16,0] = 0.5680256501163718
$vI[17,0] = 0.5687463729540241
$vI[18,0] = 0.5655328171405039
$vI[19,0] = 0.5552575731971068
$vI[20,0] = 0.5489319785686462
$vI[21,0] = 0.5522723352940062
$vI[22,0] = 0.5657073293489638
$vI[23,0] = 0.581867066373988
$ws.Range("I2:I25").Value = $vI

$vL = New-Object 'object[,]' 24,1
$vL[0,0] = 0.201913699169971
$vL[1,0] = 0.1992950126146198
$vL[2,0] = 0.1977887047952507
$vL[3,0] = 0.1972004707197854
$vL[4,0] = 0.1971043426565728
$vL[5,0] = 0.1977806679417995
$vL[6,0] = 0.2009897263537539
$vL[7,0] = 0.2080867307858796
$vL[8,0] = 0.2137893374267179
$vL[9,0] = 0.2164893477050072
$vL[10,0] = 0.2175269510048707
$vL[11,0] = 0.2173028108495032
$vL[12,0] = 0.216574408280664
$vL[13,0] = 0.2161302143540667
$vL[14,0] = 0.2136150109629682
$vL[15,0] = 0.2120990905284117
$vL[16,0] = 0.2112371410992466
$vL[17,0] = 0.2109470135190605
$vL[18,0] = 0.2122594316070945
$vL[19,0] = 0.21678794655773
$vL[20,0] = 0.2198359779356025
$vL[21,0] = 0.2182011172723861
$vL[22,0] = 0.2121869115815826
$vL[23,0] = 0.2060809034196041
$ws.Range("L2:L25").Value = $vL

$vM = New-Object 'object[,]' 24,1
$vM[0,0] = 0.2183375320802696
$vM[1,0] = 0.2011436947125844
$vM[2,0] = 0.1906086840035925
$vM[3,0] = 0.1863214146747865
$vM[4,0] = 0.1856098775143948
$vM[5,0] = 0.1905508403967957
$vM[6,0] = 0.2124046785815139
$vM[7,0] = 0.2554249684353636
$vM[8,0] = 0.2871221465328802
$vM[9,0] = 0.3015596127715341
$vM[10,0] = 0.3070290957004218
$vM[11,0] = 0.3058510458869677
$vM[12,0] = 0.3020095454721101
$vM[13,0] = 0.299656809880247
$vM[14,0] = 0.2861789660908585
$vM[15,0] = 0.27791520982489
$vM[16,0] = 0.2731638522925337
$vM[17,0] = 0.2715554316324287
$vM[18,0] = 0.2787947232307033
$vM[19,0] = 0.3031378264538134
$vM[20,0] = 0.3190608828283743
$vM[21,0] = 0.3105613170188661
$vM[22,0] = 0.2783970964873674
$vM[23,0] = 0.2437702612591437
$ws.Range("M2:M25").Value = $vM

$vO = New-Object 'object[,]' 24,1
$vO[0,0] = 2.170888144921818
$vO[1,0] = 2.187073318100701
$vO[2,0] = 2.19872333579562
$vO[3,0] = 2.203900773003937
$vO[4,0] = 2.204786433422399
$vO[5,0] = 2.198791420432059
$vO[6,0] = 2.176113051049697
$vO[7,0] = 2.14525539860972
$vO[8,0] = 2.130925810033716
$vO[9,0] = 2.126227265008851
$vO[10,0] = 2.12471046599606
$vO[11,0] = 2.125025452452832
$vO[12,0] = 2.126097213309208
$vO[13,0] = 2.126787896121868
$vO[14,0] = 2.13126955652362
$vO[15,0] = 2.134485581328164
$vO[16,0] = 2.136506623589611
$vO[17,0] = 2.137220307639865
$vO[18,0] = 2.134125499003403
$vO[19,0] = 2.125775282991782
$vO[20,0] = 2.121847952439282
$vO[21,0] = 2.123803810607257
$vO[22,0] = 2.134287756173023
$vO[23,0] = 2.152141348619523
$ws.Range("O2:O25").Value = $vO

